$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 365.3125
$ws.Range("I11").Value = 365.3125
$ws.Range("K11").Value = 365.3125
$ws.Range("M11").Value = -225.3125
$ws.Range("H17").Value = 1201.8718
$ws.Range("J17").Value = 1180.8684
$ws.Range("L17").Value = 3542.6052
$ws.Range("N17").Value = -3878.6052
$ws.Range("H18").Value = 83335690
$ws.Range("I18").Value = 83335690
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 83335690
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -83335406
$ws.Range("N18").ClearContents()
$ws.Range("H43").Value = 2443
$ws.Range("I43").Value = 2959.6667
$ws.Range("K43").Value = 2959.6667
$ws.Range("M43").Value = -2890.6667
$ws.Range("H74").Value = 3560.5454
$ws.Range("I74").Value = 3560.5454
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 3560.5454
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -2624.5454
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 3560.5454
$ws.Range("I77").Value = 3560.5454
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 17802.727
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -13122.727
$ws.Range("N77").ClearContents()
$ws.Range("H86").Value = 1962.1538
$ws.Range("I86").Value = 2126.3
$ws.Range("K86").Value = 2126.3
$ws.Range("M86").Value = -1003.3
$ws.Range("H89").Value = 1962.1538
$ws.Range("I89").Value = 2126.3
$ws.Range("K89").Value = 10631.5
$ws.Range("M89").Value = -5015.5
$ws.Range("H94").Value = 11906979
$ws.Range("I94").Value = 12989431
$ws.Range("K94").Value = 12989431
$ws.Range("M94").Value = -12988980
$ws.Range("H96").Value = 686.56525
$ws.Range("I96").Value = 459.5
$ws.Range("K96").Value = 1378.5
$ws.Range("M96").Value = -5.5
$ws.Range("H98").Value = 1999
$ws.Range("I98").Value = 1999
$ws.Range("K98").Value = 1999
$ws.Range("M98").Value = -501
$ws.Range("H106").Value = 15877218
$ws.Range("I106").Value = 19609862
$ws.Range("K106").Value = 19609862
$ws.Range("M106").Value = -19609231
$ws.Range("H107").Value = 18522416
$ws.Range("J107").Value = 5488.4443
$ws.Range("L107").Value = 5488.4443
$ws.Range("N107").Value = -9328.444299999999
$ws.Range("H112").Value = 1297.2449
$ws.Range("I112").Value = 1899.6666
$ws.Range("K112").Value = 5698.9998
$ws.Range("M112").Value = -4590.9998
$ws.Range("H122").Value = 1999
$ws.Range("I122").Value = 1999
$ws.Range("K122").Value = 5997
$ws.Range("M122").Value = -3547
$ws.Range("H132").Value = 9317.710999999999
$ws.Range("I132").Value = 5749.5
$ws.Range("J132").Value = 17048.834
$ws.Range("K132").Value = 17248.5
$ws.Range("L132").Value = 51146.50199999999
$ws.Range("M132").Value = -14718.5
$ws.Range("N132").Value = -56206.50199999999
$ws.Range("H136").Value = 118998.5
$ws.Range("I136").Value = 80779
$ws.Range("J136").Value = 131738.33
$ws.Range("K136").Value = 80779
$ws.Range("L136").Value = 131738.33
$ws.Range("M136").Value = -75679
$ws.Range("N136").Value = -141938.33
$ws.Range("H137").Value = 4231.0977
$ws.Range("I137").Value = 6301.2383
$ws.Range("J137").Value = 2057.45
$ws.Range("K137").Value = 18903.7149
$ws.Range("L137").Value = 6172.349999999999
$ws.Range("M137").Value = -16353.7149
$ws.Range("N137").Value = -11272.35
$ws.Range("H141").Value = 11800.066
$ws.Range("I141").Value = 14733.444
$ws.Range("K141").Value = 44200.33199999999
$ws.Range("M141").Value = -39020.33199999999
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 29367400
$ws.Range("I2").Value = 48053764
$ws.Range("J2").Value = 3109.1428
$ws.Range("K2").Value = 48053764
$ws.Range("L2").Value = 3109.1428
$ws.Range("M2").Value = -48053651
$ws.Range("N2").Value = -3335.1428
$ws.Range("H4").Value = 529.1111
$ws.Range("I4").Value = 229.54546
$ws.Range("J4").Value = 999.8570999999999
$ws.Range("K4").Value = 229.54546
$ws.Range("L4").Value = 999.8570999999999
$ws.Range("M4").Value = -113.54546
$ws.Range("N4").Value = -1231.8571
$ws.Range("H32").Value = 3598.178
$ws.Range("I32").Value = 3196.8
$ws.Range("K32").Value = 3196.8
$ws.Range("M32").Value = -2909.8
$ws.Range("H45").Value = 7589.769
$ws.Range("J45").Value = 4299.25
$ws.Range("L45").Value = 4299.25
$ws.Range("N45").Value = -5053.25
$ws.Range("H61").Value = 3127.0264
$ws.Range("I61").Value = 3025.7812
$ws.Range("J61").Value = 3667
$ws.Range("K61").Value = 3025.7812
$ws.Range("L61").Value = 3667
$ws.Range("M61").Value = -2813.7812
$ws.Range("N61").Value = -4091
$ws.Range("H74").Value = 3865.5374
$ws.Range("I74").Value = 3911.2969
$ws.Range("J74").Value = 2889.3333
$ws.Range("K74").Value = 3911.2969
$ws.Range("L74").Value = 2889.3333
$ws.Range("M74").Value = -3037.2969
$ws.Range("N74").Value = -4637.3333
$ws.Range("H77").Value = 3865.5374
$ws.Range("I77").Value = 3911.2969
$ws.Range("J77").Value = 2889.3333
$ws.Range("K77").Value = 19556.4845
$ws.Range("L77").Value = 14446.6665
$ws.Range("M77").Value = -15188.4845
$ws.Range("N77").Value = -23182.6665
$ws.Range("H88").Value = 1763
$ws.Range("I88").Value = 1514.75
$ws.Range("J88").Value = 1961.6
$ws.Range("K88").Value = 1514.75
$ws.Range("L88").Value = 1961.6
$ws.Range("M88").Value = -1108.75
$ws.Range("N88").Value = -2773.6
$ws.Range("H91").Value = 1763
$ws.Range("I91").Value = 1514.75
$ws.Range("J91").Value = 1961.6
$ws.Range("K91").Value = 1514.75
$ws.Range("L91").Value = 1961.6
$ws.Range("M91").Value = -110.75
$ws.Range("N91").Value = -4769.6
$ws.Range("H97").Value = 58888760
$ws.Range("I97").Value = 83339740
$ws.Range("J97").Value = 206397.8
$ws.Range("K97").Value = 83339740
$ws.Range("L97").Value = 206397.8
$ws.Range("M97").Value = -83339244
$ws.Range("N97").Value = -207389.8
$ws.Range("H110").Value = 2860.5
$ws.Range("I110").Value = 3030
$ws.Range("J110").Value = 2013
$ws.Range("K110").Value = 3030
$ws.Range("L110").Value = 2013
$ws.Range("M110").Value = -985
$ws.Range("N110").Value = -6103
$ws.Range("H116").Value = 29367400
$ws.Range("I116").Value = 48053764
$ws.Range("J116").Value = 3109.1428
$ws.Range("K116").Value = 48053764
$ws.Range("L116").Value = 3109.1428
$ws.Range("M116").Value = -48051470
$ws.Range("N116").Value = -7697.1428
$ws.Range("H122").Value = 3595.2307
$ws.Range("I122").Value = 4020
$ws.Range("J122").Value = 3329.75
$ws.Range("K122").Value = 12060
$ws.Range("L122").Value = 9989.25
$ws.Range("M122").Value = -9610
$ws.Range("N122").Value = -14889.25
$ws.Range("H132").Value = 3940.0513
$ws.Range("I132").Value = 3215.5
$ws.Range("J132").Value = 6355.222
$ws.Range("K132").Value = 9646.5
$ws.Range("L132").Value = 19065.666
$ws.Range("M132").Value = -7116.5
$ws.Range("N132").Value = -24125.666
$ws.Range("H135").Value = 150000
$ws.Range("J135").Value = 150000
$ws.Range("L135").Value = 150000
$ws.Range("N135").Value = -160140
$ws.Range("H136").Value = 3127.0264
$ws.Range("I136").Value = 3025.7812
$ws.Range("J136").Value = 3667
$ws.Range("K136").Value = 9077.3436
$ws.Range("L136").Value = 11001
$ws.Range("M136").Value = -6527.3436
$ws.Range("N136").Value = -16101
$ws.Range("H138").Value = 69557.5
$ws.Range("J138").Value = 69557.5
$ws.Range("L138").Value = 69557.5
$ws.Range("N138").Value = -79837.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 29367400
$ws.Range("I3").Value = 48053764
$ws.Range("J3").Value = 3109.1428
$ws.Range("K3").Value = 48053764
$ws.Range("L3").Value = 3109.1428
$ws.Range("M3").Value = -48053650
$ws.Range("N3").Value = -3337.1428
$ws.Range("H20").Value = 3434.8333
$ws.Range("I20").Value = 3389.4285
$ws.Range("K20").Value = 3389.4285
$ws.Range("M20").Value = -3142.4285
$ws.Range("H22").Value = 28044.666
$ws.Range("I22").Value = 210.38461
$ws.Range("J22").Value = 100413.8
$ws.Range("K22").Value = 210.38461
$ws.Range("L22").Value = 100413.8
$ws.Range("M22").Value = -37.38461000000001
$ws.Range("N22").Value = -100759.8
$ws.Range("H107").Value = 2604.647
$ws.Range("I107").Value = 3783.5
$ws.Range("J107").Value = 1961.6364
$ws.Range("K107").Value = 3783.5
$ws.Range("L107").Value = 1961.6364
$ws.Range("M107").Value = -1863.5
$ws.Range("N107").Value = -5801.6364
$ws.Range("H115").Value = 100000
$ws.Range("J115").Value = 100000
$ws.Range("L115").Value = 100000
$ws.Range("N115").Value = -103134
$ws.Range("H134").Value = 3196.7222
$ws.Range("I134").Value = 3214.1765
$ws.Range("K134").Value = 9642.529500000001
$ws.Range("M134").Value = -7107.529500000001
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1491.4667
$ws.Range("I16").Value = 935.4
$ws.Range("J16").Value = 1769.5
$ws.Range("K16").Value = 935.4
$ws.Range("L16").Value = 1769.5
$ws.Range("M16").Value = -648.4
$ws.Range("N16").Value = -2343.5
$ws.Range("H20").Value = 133074.67
$ws.Range("J20").Value = 133074.67
$ws.Range("L20").Value = 133074.67
$ws.Range("N20").Value = -133546.67
$ws.Range("H30").Value = 133074.67
$ws.Range("J30").Value = 133074.67
$ws.Range("L30").Value = 133074.67
$ws.Range("N30").Value = -133256.67
$ws.Range("H31").Value = 1866.591
$ws.Range("I31").Value = 1398.5
$ws.Range("J31").Value = 2685.75
$ws.Range("K31").Value = 1398.5
$ws.Range("L31").Value = 2685.75
$ws.Range("M31").Value = -1103.5
$ws.Range("N31").Value = -3275.75
$ws.Range("H34").Value = 1866.591
$ws.Range("I34").Value = 1398.5
$ws.Range("J34").Value = 2685.75
$ws.Range("K34").Value = 1398.5
$ws.Range("L34").Value = 2685.75
$ws.Range("M34").Value = -1196.5
$ws.Range("N34").Value = -3089.75
$ws.Range("H56").Value = 70998.25
$ws.Range("J56").Value = 70000
$ws.Range("L56").Value = 70000
$ws.Range("N56").Value = -71690
$ws.Range("H62").Value = 83336990
$ws.Range("J62").Value = 4507.8335
$ws.Range("L62").Value = 4507.8335
$ws.Range("N62").Value = -5755.8335
$ws.Range("H65").Value = 83336990
$ws.Range("J65").Value = 4507.8335
$ws.Range("L65").Value = 22539.1675
$ws.Range("N65").Value = -28779.1675
$ws.Range("H75").Value = 38662.668
$ws.Range("J75").Value = 38662.668
$ws.Range("L75").Value = 38662.668
$ws.Range("N75").Value = -40658.668
$ws.Range("H78").Value = 38662.668
$ws.Range("J78").Value = 38662.668
$ws.Range("L78").Value = 115988.004
$ws.Range("N78").Value = -125972.004
$ws.Range("H99").Value = 6696.0684
$ws.Range("I99").Value = 7957.3335
$ws.Range("J99").Value = 6371.7427
$ws.Range("K99").Value = 7957.3335
$ws.Range("L99").Value = 6371.7427
$ws.Range("M99").Value = -6459.3335
$ws.Range("N99").Value = -9367.742699999999
$ws.Range("H107").Value = 1720.8889
$ws.Range("J107").Value = 2133.5
$ws.Range("L107").Value = 2133.5
$ws.Range("N107").Value = -5973.5
$ws.Range("H110").Value = 139998.2
$ws.Range("J110").Value = 139998.2
$ws.Range("L110").Value = 139998.2
$ws.Range("N110").Value = -148178.2
$ws.Range("H113").Value = 1491.4667
$ws.Range("I113").Value = 935.4
$ws.Range("J113").Value = 1769.5
$ws.Range("K113").Value = 935.4
$ws.Range("L113").Value = 1769.5
$ws.Range("M113").Value = 1234.6
$ws.Range("N113").Value = -6109.5
$ws.Range("H122").Value = 4087
$ws.Range("I122").Value = 4233.9
$ws.Range("J122").Value = 3903.375
$ws.Range("K122").Value = 12701.7
$ws.Range("L122").Value = 11710.125
$ws.Range("M122").Value = -10251.7
$ws.Range("N122").Value = -16610.125
$ws.Range("H126").Value = 6696.0684
$ws.Range("I126").Value = 7957.3335
$ws.Range("J126").Value = 6371.7427
$ws.Range("K126").Value = 23872.0005
$ws.Range("L126").Value = 19115.2281
$ws.Range("M126").Value = -21402.0005
$ws.Range("N126").Value = -24055.2281
$ws.Range("H128").Value = 133074.67
$ws.Range("J128").Value = 133074.67
$ws.Range("L128").Value = 133074.67
$ws.Range("N128").Value = -143034.67
$ws.Range("H133").Value = 154996
$ws.Range("J133").Value = 154996
$ws.Range("L133").Value = 154996
$ws.Range("N133").Value = -160056
$ws.Range("H134").Value = 7000
$ws.Range("I134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("M134").ClearContents()
$ws.Range("H138").Value = 99473.336
$ws.Range("J138").Value = 99473.336
$ws.Range("L138").Value = 99473.336
$ws.Range("N138").Value = -109753.336
$ws.Range("H140").Value = 107265.555
$ws.Range("J140").Value = 109912.836
$ws.Range("L140").Value = 109912.836
$ws.Range("N140").Value = -120272.836
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 2895.0625
$ws.Range("J11").Value = 3015.9285
$ws.Range("L11").Value = 9047.7855
$ws.Range("N11").Value = -9327.7855
$ws.Range("H37").Value = 82112.45
$ws.Range("J37").Value = 82112.45
$ws.Range("L37").Value = 246337.35
$ws.Range("N37").Value = -246561.35
$ws.Range("H57").Value = 5859.7
$ws.Range("I57").Value = 8000
$ws.Range("J57").Value = 5324.625
$ws.Range("K57").Value = 24000
$ws.Range("L57").Value = 15973.875
$ws.Range("M57").Value = -23441
$ws.Range("N57").Value = -17091.875
$ws.Range("H92").Value = 665.9
$ws.Range("J92").Value = 608.1667
$ws.Range("L92").Value = 1824.5001
$ws.Range("N92").Value = -4320.5001
$ws.Range("H97").Value = 570
$ws.Range("J97").Value = 536.1667
$ws.Range("L97").Value = 1608.5001
$ws.Range("N97").Value = -2600.5001
$ws.Range("H132").Value = 3342.4443
$ws.Range("I132").Value = 3023
$ws.Range("K132").Value = 27207
$ws.Range("M132").Value = -24677
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 60.42105
$ws.Range("I2").Value = 44.214287
$ws.Range("K2").Value = 44.214287
$ws.Range("M2").Value = 68.785713
$ws.Range("H70").Value = 7874
$ws.Range("I70").Value = 7874
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 7874
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -7604
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 7874
$ws.Range("I73").Value = 7874
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 7874
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -6938
$ws.Range("N73").ClearContents()
$ws.Range("H122").Value = 8051
$ws.Range("I122").Value = 9671
$ws.Range("J122").Value = 7079
$ws.Range("K122").Value = 29013
$ws.Range("L122").Value = 21237
$ws.Range("M122").Value = -26563
$ws.Range("N122").Value = -26137
$ws.Range("H132").Value = 5636.2163
$ws.Range("I132").Value = 5187.696
$ws.Range("J132").Value = 6373.0713
$ws.Range("K132").Value = 15563.088
$ws.Range("L132").Value = 19119.2139
$ws.Range("M132").Value = -13033.088
$ws.Range("N132").Value = -24179.2139
$ws.Range("H141").Value = 42041
$ws.Range("J141").Value = 42041
$ws.Range("L141").Value = 42041
$ws.Range("N141").Value = -52401
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2142
$ws.Range("I7").Value = 1977.5
$ws.Range("J7").Value = 2800
$ws.Range("K7").Value = 1977.5
$ws.Range("L7").Value = 2800
$ws.Range("M7").Value = -1865.5
$ws.Range("N7").Value = -3024
$ws.Range("H40").Value = 4527.6665
$ws.Range("I40").Value = 4374.577
$ws.Range("K40").Value = 4374.577
$ws.Range("M40").Value = -4238.577
$ws.Range("H82").Value = 1975.79
$ws.Range("I82").Value = 1948.2188
$ws.Range("J82").Value = 2637.5
$ws.Range("K82").Value = 1948.2188
$ws.Range("L82").Value = 2637.5
$ws.Range("M82").Value = -1587.2188
$ws.Range("N82").Value = -3359.5
$ws.Range("H85").Value = 1975.79
$ws.Range("I85").Value = 1948.2188
$ws.Range("J85").Value = 2637.5
$ws.Range("K85").Value = 1948.2188
$ws.Range("L85").Value = 2637.5
$ws.Range("M85").Value = -700.2188000000001
$ws.Range("N85").Value = -5133.5
$ws.Range("H126").Value = 2142
$ws.Range("I126").Value = 1977.5
$ws.Range("J126").Value = 2800
$ws.Range("K126").Value = 5932.5
$ws.Range("L126").Value = 8400
$ws.Range("M126").Value = -3462.5
$ws.Range("N126").Value = -13340
$ws.Range("H132").Value = 18262.42
$ws.Range("I132").Value = 22042.814
$ws.Range("J132").Value = 4653
$ws.Range("K132").Value = 66128.442
$ws.Range("L132").Value = 13959
$ws.Range("M132").Value = -63598.442
$ws.Range("N132").Value = -19019
$ws.Range("H136").Value = 3168284
$ws.Range("I136").Value = 4296091
$ws.Range("J136").Value = 10424.467
$ws.Range("K136").Value = 12888273
$ws.Range("L136").Value = 31273.401
$ws.Range("M136").Value = -12885723
$ws.Range("N136").Value = -36373.401
$ws.Range("H139").Value = 119999
$ws.Range("J139").Value = 119999
$ws.Range("L139").Value = 119999
$ws.Range("N139").Value = -130279
$ws.Range("H140").Value = 97394.8
$ws.Range("J140").Value = 97394.8
$ws.Range("L140").Value = 97394.8
$ws.Range("N140").Value = -107754.8
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 97433.164
$ws.Range("J46").Value = 97433.164
$ws.Range("L46").Value = 97433.164
$ws.Range("N46").Value = -97895.164
$ws.Range("H49").Value = 11111
$ws.Range("I49").Value = 11111
$ws.Range("K49").Value = 11111
$ws.Range("M49").Value = -10881
$ws.Range("H107").Value = 4196.7856
$ws.Range("I107").Value = 3987.111
$ws.Range("K107").Value = 11961.333
$ws.Range("M107").Value = -10041.333
$ws.Range("H122").Value = 11737.333
$ws.Range("J122").Value = 33665.668
$ws.Range("L122").Value = 100997.004
$ws.Range("N122").Value = -105897.004
$ws.Range("H126").Value = 5758.1816
$ws.Range("I126").Value = 4079.875
$ws.Range("J126").Value = 10233.667
$ws.Range("K126").Value = 12239.625
$ws.Range("L126").Value = 30701.001
$ws.Range("M126").Value = -9769.625
$ws.Range("N126").Value = -35641.001
$ws.Range("H132").Value = 2526.4285
$ws.Range("I132").Value = 2039.0217
$ws.Range("J132").Value = 10000
$ws.Range("K132").Value = 6117.0651
$ws.Range("L132").Value = 30000
$ws.Range("M132").Value = -3587.0651
$ws.Range("N132").Value = -35060
$ws.Range("H134").Value = 97433.164
$ws.Range("J134").Value = 97433.164
$ws.Range("L134").Value = 292299.492
$ws.Range("N134").Value = -297369.492
$ws.Range("H136").Value = 1719.5873
$ws.Range("I136").Value = 1647.3966
$ws.Range("J136").Value = 2557
$ws.Range("K136").Value = 4942.1898
$ws.Range("L136").Value = 7671
$ws.Range("M136").Value = -2392.1898
$ws.Range("N136").Value = -12771
$ws.Range("H137").Value = 103294.336
$ws.Range("J137").Value = 103294.336
$ws.Range("L137").Value = 103294.336
$ws.Range("N137").Value = -113494.336
$ws.Range("H140").Value = 204896.33
$ws.Range("J140").Value = 204896.33
$ws.Range("L140").Value = 204896.33
$ws.Range("N140").Value = -215256.33
